# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) / "Valor Mora" (column F) table in
# B16:J54 is re-ordered so the periods run in ascending order
# (1607 .. 1909) instead of the original descending order
# (1909 .. 1607). The Valor Mora amount stays attached to its period,
# it just slides to the new row. Row-level formatting (borders, etc.)
# must stay exactly where it was, so the cells are rewritten in place
# rather than using a Range.Sort (which would drag the per-row style
# along with the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Periodo Mora" order (ascending) for B16:J54.
$periodos = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909"
)

# Matching "Valor Mora" amounts - the same value that used to travel
# with each period, now following it to its new row.
$valores = @(
    24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,27083
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
